$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Term"
$ws.Range("B2").Value = "Rchn & Rcssd"
$ws.Range("J2").Value = '["FFS/PPO/ACO/HMO/Medi-Cal"]'
$ws.Range("K2").Value = "Information not found"
$ws.Range("S2").Value = "prajay.sapkale@hilabs.com"

# Row 3
$ws.Range("A3").Value = "Term"
$ws.Range("B3").Value = "Cyrus Hendricks, M.D."
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4568888895"
$ws.Range("G3").Value = "207R00000X"
$ws.Range("J3").Value = "PPG#’s, Medicare"
$ws.Range("L3").Value = "Mercian Medical Group – 0P4"
$ws.Range("S3").Value = "prajay.sapkale@hilabs.com"

# Row 4
$ws.Range("A4").Value = "Term"
$ws.Range("B4").Value = "Paul Mcmallan, Md"
$ws.Range("S4").Value = "p rajay.sapkale@hilabs.com"
